$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 83
$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

# "25" looks numeric, so Excel would otherwise store it as a number.
# Force it to be stored as text (matching the other rows), then restore
# the default (Normal) cell style so no stray formatting is left behind.
$ws.Cells.Item($row, 3).Value = "'25"
$ws.Cells.Item($row, 3).Style = "Normal"

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
